# Fixed way of SMS app validation in intent
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- G9: add ScrollPage(runtest_top_xpath); before final TakeScreenshot ---
$ws.Range("G9").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0883_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nSwitchApp(NATIVE_APP);`nDrawSignature(inlinesignature_view_xpath);`nSwitchApp(WEBVIEW);`nTakeScreenshot(VT200-0883-01);`nwait(4);`nScrollPage(runtest_top_xpath);`nTakeScreenshot(VT200-0883-02);`nwait(2);`nvalidate4;"

# --- H15: remove trailing validate_Result lines ---
$ws.Range("H15").Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Signature JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0889`n};`nvalidate4`n{`nvalidate_Screenshot=VT200-0889-01`nvalidate_Screenshot=VT200-0889-02`n};"

# --- G16: remove redundant wait(2); lines ---
$ws.Range("G16").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0890_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nTakeScreenshot(VT200-0890-01);`npress_Key(Home);`nvalidate4;`nlaunch_App_Device(com.rhomobile.compliancetest_js/com.rhomobile.rhodes.RhodesActivity);`nTakeScreenshot(VT200-0890-02);`nvalidate5;"

# --- G18: fix ScrollPage target from results_xpath to runtest_top_xpath ---
$ws.Range("G18").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0892_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nScrollPage(runtest_top_xpath);`nTakeScreenshot(VT200-0892);`nwait(2);`nvalidate4;"

# --- Row 15 height shrinks because H15 text got shorter ---
$ws.Rows.Item(15).RowHeight = 192

# --- Selection moved from D2 to C2 ---
$ws.Range("C2").Select()
